# Commit: "#5: property aircraft done"
#
# The "建物" (building) property sheet had its "property_category" column
# (column I) erroneously populated with the string "land" (copy/paste
# leftover from the "土地" land sheet) for every data row. This fixes it
# so the building sheet correctly reports "building" as its category.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("建物")

# Data rows 2-10, column I ("property_category") -> "building"
$ws.Range("I2:I10").Value = "building"
